$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999446620147,
    0.99414982016468223,
    0.97172387839909891,
    0.96171981562347542,
    0.95211972498379849,
    0.92861073718331111,
    0.92576710356951053,
    0.92249619870647748,
    0.92073180483135419,
    0.91227946431830587,
    0.91107823107343244,
    0.90904892293366069,
    0.90112985610271024,
    0.89838283432854826,
    0.89717194571699876,
    0.89632793586118142,
    0.89261992890650865,
    0.89151101007909705,
    0.99094109109290884,
    0.96472685284108772,
    0.96332833677144802,
    0.96206382183136463,
    0.96468815346922177,
    0.95166633239031162,
    0.94520917807485505,
    0.92556236947845139,
    0.921561941527409,
    0.90393019031405819,
    0.89150845395625999,
    0.88638654713653353,
    0.88232842370363418,
    0.88064909096927901,
    0.88012908342223928
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
